$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 becomes a second task row, duplicated (format-wise) from row 12 ---
# Copy the formatting of row 12's cells onto row 13 (matches the author having
# inserted a new task row by copying the template row and editing its text).
$cols = @("A","B","C","D","E","F","G","H","I","J","K")
foreach ($col in $cols) {
    $ws.Range("$col`12").Copy()
    $ws.Range("$col`13").PasteSpecial(-4122)  # xlPasteFormats
}

# --- Task #1 (row 12): description changes to the DB objects task ---
$ws.Range("B12").Value = "Realizar la creacion de los objetos del DB TMS ""ObjetosDB.sql"" "

# --- Task #2 (row 13): new task number + merge-request description (bold runs) ---
$ws.Range("A13").Value = 2

$mergeText = "Realizar el pase a producción del siguiente merge:" + "`n" + " De calidad a master, http://gitlab.estrellaroja.com.mx/java/facturacion-api/-/merge_requests/229"
$ws.Range("B13").Value = $mergeText

$boldStart1 = $mergeText.IndexOf("calidad") + 1
$ws.Range("B13").Characters($boldStart1, 7).Font.Bold = $true

$boldStart2 = $mergeText.IndexOf("master") + 1
$ws.Range("B13").Characters($boldStart2, 6).Font.Bold = $true

$ws.Range("G13").Value = "TI"
$ws.Range("H13").Value = "Pendiente"

# --- Selection / scroll position moves to the newly edited cell ---
[void]$ws.Range("B12:E12").Select()
